$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 8.598832290884035
$ws.Range("D2").Value = 4.723202461059867
$ws.Range("E2").Value = 13.38954024654233
$ws.Range("F2").Value = 23.50927051105727
$ws.Range("G2").Value = 3.624690452597635
$ws.Range("I2").Value = 21.09569803924334
$ws.Range("L2").Value = 9.9242591202757
$ws.Range("M2").Value = 59.2029173905363
$ws.Range("O2").Value = 21.02348849842801
# Row 3
$ws.Range("C3").Value = 8.647245823064207
$ws.Range("D3").Value = 4.736437437554335
$ws.Range("E3").Value = 13.12136397856194
$ws.Range("F3").Value = 23.77624624529844
$ws.Range("G3").Value = 3.628633911525839
$ws.Range("I3").Value = 21.47565035685377
$ws.Range("L3").Value = 9.77756997543413
$ws.Range("M3").Value = 55.75968231933449
$ws.Range("O3").Value = 21.32175786957491
# Row 4
$ws.Range("C4").Value = 8.6814115622554
$ws.Range("D4").Value = 4.745702001242872
$ws.Range("E4").Value = 12.96042131958315
$ws.Range("F4").Value = 23.95618706659218
$ws.Range("G4").Value = 3.631156404444084
$ws.Range("I4").Value = 21.72052774235683
$ws.Range("L4").Value = 9.690532490804671
$ws.Range("M4").Value = 53.52543047524149
$ws.Range("O4").Value = 21.5169867784298
# Row 5
$ws.Range("C5").Value = 8.696430274318075
$ws.Range("D5").Value = 4.749761735344268
$ws.Range("E5").Value = 12.89585861839808
$ws.Range("F5").Value = 24.03345047336537
$ws.Range("G5").Value = 3.632209948289412
$ws.Range("I5").Value = 21.8232051755318
$ws.Range("L5").Value = 9.655863289414098
$ws.Range("M5").Value = 52.58497548028853
$ws.Range("O5").Value = 21.5995350078382
# Row 6
$ws.Range("C6").Value = 8.698989720199791
$ws.Range("D6").Value = 4.750452972102018
$ws.Range("E6").Value = 12.88520235558879
$ws.Range("F6").Value = 24.04651510947017
$ws.Range("G6").Value = 3.632386439724482
$ws.Range("I6").Value = 21.84042841221877
$ws.Range("L6").Value = 9.650155836312866
$ws.Range("M6").Value = 52.42701088372174
$ws.Range("O6").Value = 21.6134213278094
# Row 7
$ws.Range("C7").Value = 8.681609699289867
$ws.Range("D7").Value = 4.745755603315625
$ws.Range("E7").Value = 12.95954634779335
$ws.Range("F7").Value = 23.95721324972958
$ws.Range("G7").Value = 3.631170509009972
$ws.Range("I7").Value = 21.72190082471742
$ws.Range("L7").Value = 9.690061645131024
$ws.Range("M7").Value = 53.51286821074513
$ws.Range("O7").Value = 21.51808800597797
# Row 8
$ws.Range("C8").Value = 8.614591352473361
$ws.Range("D8").Value = 4.72752848774681
$ws.Range("E8").Value = 13.29635384934684
$ws.Range("F8").Value = 23.59794366789108
$ws.Range("G8").Value = 3.626029253621387
$ws.Range("I8").Value = 21.22428522521629
$ws.Range("L8").Value = 9.87307446679405
$ws.Range("M8").Value = 58.04072351738681
$ws.Range("O8").Value = 21.12379193168398
# Row 9
$ws.Range("C9").Value = 8.51933797884606
$ws.Range("D9").Value = 4.700905714455584
$ws.Range("E9").Value = 13.98251699479431
$ws.Range("F9").Value = 23.02472523012903
$ws.Range("G9").Value = 3.616742536032985
$ws.Range("I9").Value = 20.34163218955256
$ws.Range("L9").Value = 10.2543298079292
$ws.Range("M9").Value = 65.96086883468409
$ws.Range("O9").Value = 20.44885765914895
# Row 10
$ws.Range("C10").Value = 8.472733712605933
$ws.Range("D10").Value = 4.687036825225635
$ws.Range("E10").Value = 14.49731929680247
$ws.Range("F10").Value = 22.68960973323085
$ws.Range("G10").Value = 3.61039326998938
$ws.Range("I10").Value = 19.75186723456067
$ws.Range("L10").Value = 10.54581748743339
$ws.Range("M10").Value = 71.19099203355204
$ws.Range("O10").Value = 20.01635005597337
# Row 11
$ws.Range("C11").Value = 8.456903566672493
$ws.Range("D11").Value = 4.681993653723731
$ws.Range("E11").Value = 14.73284239939688
$ws.Range("F11").Value = 22.55717890065291
$ws.Range("G11").Value = 3.607605170485251
$ws.Range("I11").Value = 19.49681294644806
$ws.Range("L11").Value = 10.680421551188
$ws.Range("M11").Value = 73.44243390862134
$ws.Range("O11").Value = 19.83417378207954
# Row 12
$ws.Range("C12").Value = 8.451706179274792
$ws.Range("D12").Value = 4.680268617268696
$ws.Range("E12").Value = 14.82214400269012
$ws.Range("F12").Value = 22.51002334945662
$ws.Range("G12").Value = 3.606563595360479
$ws.Range("I12").Value = 19.40218172741332
$ws.Range("L12").Value = 10.73164256781532
$ws.Range("M12").Value = 74.27665689385766
$ws.Range("O12").Value = 19.76735860645756
# Row 13
$ws.Range("C13").Value = 8.45278967604699
$ws.Range("D13").Value = 4.680631876087644
$ws.Range("E13").Value = 14.8029074177746
$ws.Range("F13").Value = 22.52004416374176
$ws.Range("G13").Value = 3.606787288005969
$ws.Range("I13").Value = 19.42247457655541
$ws.Range("L13").Value = 10.72060071777356
$ws.Range("M13").Value = 74.09780749544129
$ws.Range("O13").Value = 19.78165064076838
# Row 14
$ws.Range("C14").Value = 8.456459854677908
$ws.Range("D14").Value = 4.681848016480967
$ws.Range("E14").Value = 14.74018743812045
$ws.Range("F14").Value = 22.55323873097222
$ws.Range("G14").Value = 3.607519195659128
$ws.Range("I14").Value = 19.48898814885783
$ws.Range("L14").Value = 10.68463075477188
$ws.Range("M14").Value = 73.51143345738107
$ws.Range("O14").Value = 19.82863289782007
# Row 15
$ws.Range("C15").Value = 8.458812502023635
$ws.Range("D15").Value = 4.682617072852871
$ws.Range("E15").Value = 14.70178227028006
$ws.Range("F15").Value = 22.57396461088776
$ws.Range("G15").Value = 3.607969355952104
$ws.Range("I15").Value = 19.52998544332287
$ws.Range("L15").Value = 10.66262946058757
$ws.Range("M15").Value = 73.14987394766126
$ws.Range("O15").Value = 19.85769593834329
# Row 16
$ws.Range("C16").Value = 8.473878434526506
$ws.Range("D16").Value = 4.687392106392401
$ws.Range("E16").Value = 14.48194741499947
$ws.Range("F16").Value = 22.69867692532949
$ws.Range("G16").Value = 3.610577479051358
$ws.Range("I16").Value = 19.76880601017704
$ws.Range("L16").Value = 10.53705777863067
$ws.Range("M16").Value = 71.04128544166348
$ws.Range("O16").Value = 20.02855529331172
# Row 17
$ws.Range("C17").Value = 8.484513963602575
$ws.Range("D17").Value = 4.690647514602539
$ws.Range("E17").Value = 14.34737108106834
$ws.Range("F17").Value = 22.78039121819747
$ws.Range("G17").Value = 3.612203009863404
$ws.Range("I17").Value = 19.91873737748604
$ws.Range("L17").Value = 10.46050897918525
$ws.Range("M17").Value = 69.71501894571077
$ws.Range("O17").Value = 20.13715671544673
# Row 18
$ws.Range("C18").Value = 8.491135208425247
$ws.Range("D18").Value = 4.692638920104686
$ws.Range("E18").Value = 14.2700965443335
$ws.Range("F18").Value = 22.82926796882584
$ws.Range("G18").Value = 3.613147416677591
$ws.Range("I18").Value = 20.0062161844882
$ws.Range("L18").Value = 10.41667109196836
$ws.Range("M18").Value = 68.9401562486316
$ws.Range("O18").Value = 20.20098831771997
# Row 19
$ws.Range("C19").Value = 8.493462853318233
$ws.Range("D19").Value = 4.693333531386934
$ws.Range("E19").Value = 14.24395754187696
$ws.Range("F19").Value = 22.84613568416481
$ws.Range("G19").Value = 3.613468804473424
$ws.Range("I19").Value = 20.03604694572618
$ws.Range("L19").Value = 10.40186234268132
$ws.Range("M19").Value = 68.67573471397267
$ws.Range("O19").Value = 20.22283315692814
# Row 20
$ws.Range("C20").Value = 8.483329455332402
$ws.Range("D20").Value = 4.690288637804862
$ws.Range("E20").Value = 14.3616840675487
$ws.Range("F20").Value = 22.77149752010259
$ws.Range("G20").Value = 3.612028993243182
$ws.Range("I20").Value = 19.90264801874758
$ws.Range("L20").Value = 10.46863825270995
$ws.Range("M20").Value = 69.8574471725571
$ws.Range("O20").Value = 20.12545390551846
# Row 21
$ws.Range("C21").Value = 8.455359999635444
$ws.Range("D21").Value = 4.681485771813779
$ws.Range("E21").Value = 14.75860731623032
$ws.Range("F21").Value = 22.54340652046113
$ws.Range("G21").Value = 3.607303832344465
$ws.Range("I21").Value = 19.46939809317758
$ws.Range("L21").Value = 10.69518953941037
$ws.Range("M21").Value = 73.68416340749565
$ws.Range("O21").Value = 19.81477352405082
# Row 22
$ws.Range("C22").Value = 8.441738520748091
$ws.Range("D22").Value = 4.676810512453026
$ws.Range("E22").Value = 15.0186523011972
$ws.Range("F22").Value = 22.41183713609223
$ws.Range("G22").Value = 3.604298443807144
$ws.Range("I22").Value = 19.19765084177758
$ws.Range("L22").Value = 10.84469055617642
$ws.Range("M22").Value = 76.07821809143161
$ws.Range("O22").Value = 19.62441883899899
# Row 23
$ws.Range("C23").Value = 8.448574017714696
$ws.Range("D23").Value = 4.679206239908709
$ws.Range("E23").Value = 14.87982795053314
$ws.Range("F23").Value = 22.4804182132987
$ws.Range("G23").Value = 3.605894967498949
$ws.Range("I23").Value = 19.34162592336161
$ws.Range("L23").Value = 10.76478004832917
$ws.Range("M23").Value = 74.81023574756843
$ws.Range("O23").Value = 19.7248271324891
# Row 24
$ws.Range("C24").Value = 8.483863395028143
$ws.Range("D24").Value = 4.690450512968576
$ws.Range("E24").Value = 14.35521286387502
$ws.Range("F24").Value = 22.77551245748891
$ws.Range("G24").Value = 3.612107635407777
$ws.Range("I24").Value = 19.90991803064336
$ws.Range("L24").Value = 10.46496247243588
$ws.Range("M24").Value = 69.79309391628506
$ws.Range("O24").Value = 20.13074040580044
# Row 25
$ws.Range("C25").Value = 8.54109907544151
$ws.Range("D25").Value = 4.707118365511213
$ws.Range("E25").Value = 13.79466623631892
$ws.Range("F25").Value = 23.1651324809447
$ws.Range("G25").Value = 3.619170830643663
$ws.Range("I25").Value = 20.57025419323331
$ws.Range("L25").Value = 10.14902879895613
$ws.Range("M25").Value = 63.92153350033747
$ws.Range("O25").Value = 20.62058161582141
